# Weekly refresh of the "Fruta / hortaliza" consolidated data:
# data is re-sorted by date (column D) and a new row (with the most
# recent date) is appended at the bottom of the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that stay constant across every data row.
$A = 11
$B = "Vega Monumental Concepción"
$C = "Bíobío"
$E = 8
$F = "Fruta"
$G = 100104
$H = "Frutos de pepita"
$I = 100104003
$J = "Membrillo"
$K = "Champion"

# New (post-edit) data for rows 2..12, in final sheet order.
$rows = @(
    @{ D = 44425; L = "Primera"; M = 100; N = 12000; O = 13000; P = 12500; Q = "`$/bandeja 18 kilos granel";  R = "Región de O'Higgins"; S = 694; T = 18 },
    @{ D = 44299; L = "Primera"; M = 100; N = 10000; O = 11000; P = 10500; Q = "`$/caja 18 kilos granel";     R = "Región del Maule";    S = 583; T = 18 },
    @{ D = 44299; L = "Segunda"; M = 50;  N = 9000;  O = 9000;  P = 9000;  Q = "`$/caja 18 kilos granel";     R = "Región del Maule";    S = 500; T = 18 },
    @{ D = 44272; L = "Primera"; M = 100; N = 9000;  O = 10000; P = 9500;  Q = "`$/caja 15 kilos granel";     R = "Región de O'Higgins"; S = 633; T = 15 },
    @{ D = 44272; L = "Segunda"; M = 50;  N = 8000;  O = 8000;  P = 8000;  Q = "`$/caja 15 kilos granel";     R = "Región de O'Higgins"; S = 533; T = 15 },
    @{ D = 44307; L = "Primera"; M = 50;  N = 10000; O = 10000; P = 10000; Q = "`$/bandeja 18 kilos granel";  R = "Región de O'Higgins"; S = 556; T = 18 },
    @{ D = 44307; L = "Segunda"; M = 50;  N = 8000;  O = 8000;  P = 8000;  Q = "`$/bandeja 18 kilos granel";  R = "Región de O'Higgins"; S = 444; T = 18 },
    @{ D = 44698; L = "Primera"; M = 50;  N = 10000; O = 10000; P = 10000; Q = "`$/caja 18 kilos granel";     R = "Región de O'Higgins"; S = 556; T = 18 },
    @{ D = 44358; L = "Primera"; M = 100; N = 11000; O = 12000; P = 11500; Q = "`$/caja 18 kilos granel";     R = "Región de O'Higgins"; S = 639; T = 18 },
    @{ D = 44363; L = "Primera"; M = 100; N = 9000;  O = 10000; P = 9500;  Q = "`$/caja 15 kilos empedrada";  R = "Región de O'Higgins"; S = 633; T = 15 },
    @{ D = 44316; L = "Primera"; M = 100; N = 9000;  O = 10000; P = 9500;  Q = "`$/caja 18 kilos granel";     R = "Región de O'Higgins"; S = 528; T = 18 }
)

$r = 2
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $A
    $ws.Cells.Item($r, 2).Value = $B
    $ws.Cells.Item($r, 3).Value = $C
    $ws.Cells.Item($r, 4).Value = $row.D
    $ws.Cells.Item($r, 5).Value = $E
    $ws.Cells.Item($r, 6).Value = $F
    $ws.Cells.Item($r, 7).Value = $G
    $ws.Cells.Item($r, 8).Value = $H
    $ws.Cells.Item($r, 9).Value = $I
    $ws.Cells.Item($r, 10).Value = $J
    $ws.Cells.Item($r, 11).Value = $K
    $ws.Cells.Item($r, 12).Value = $row.L
    $ws.Cells.Item($r, 13).Value = $row.M
    $ws.Cells.Item($r, 14).Value = $row.N
    $ws.Cells.Item($r, 15).Value = $row.O
    $ws.Cells.Item($r, 16).Value = $row.P
    $ws.Cells.Item($r, 17).Value = $row.Q
    $ws.Cells.Item($r, 18).Value = $row.R
    $ws.Cells.Item($r, 19).Value = $row.S
    $ws.Cells.Item($r, 20).Value = $row.T

    # column D uses the date-number-format style (style index 2 in the
    # original workbook) for every data row, including the newly added one.
    $ws.Cells.Item($r, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"

    $r = $r + 1
}
